$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.445.61"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.537.10"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "195.66"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "584.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.62%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.206"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.628"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.67"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("E12").Value = "  -5.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.20"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.102.04"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "663.51"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +11.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.520.26"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.548.46"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.45"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.960"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.38"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "104.30"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.37"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.92"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.52"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.14"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.32"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -9.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.74"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.72"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.57%  "
$ws.Range("E33").Value = "  -5.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.82"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.776.92"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.78"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.17%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0805"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -10.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "499.42"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.84%  "
$ws.Range("E40").Value = "  -7.35%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.134"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.369"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "34.59"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0447"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("E47").Value = "  -3.07%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.31"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("E50").Value = "  +19.74%  "
$ws.Range("E51").Value = "  +61.42%  "
